# Añade la fila 13 con los resultados SmartScore de Harsevak Sandhu Singh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A13").Value = "Harsevak Sandhu Singh_20251202_122813"
$ws.Range("B13").Value = ""
$ws.Range("C13").Value = "Harsevak Sandhu Singh"
$ws.Range("D13").Value = 18
$ws.Range("E13").Value = "Male"
$ws.Range("F13").Value = "2025-12-02 12:28:13"
$ws.Range("G13").Value = "{`n  ""portion"": 0.8,`n  ""diet"": 0.7142857142857143,`n  ""salt"": 0.4,`n  ""fat"": 0.6,`n  ""natural"": 0.6,`n  ""convenience"": 1.0,`n  ""price"": 0.6`n}"
$ws.Range("H13").Value = "Nongshim Neoguri Spicy Seafood"
$ws.Range("I13").Value = "'0.582"
$ws.Range("J13").Value = "Sabor a marisco, umami, picante equilibrado, buena textura, algo salado"
$ws.Range("K13").Value = "Nissin Chow Mein Teriyaki Beef"
$ws.Range("L13").Value = "'0.507"
$ws.Range("M13").Value = "Fácil de preparar, porción generosa, salsa suave, necesita mejoras, alto en grasa"
$ws.Range("N13").Value = "Nongshim Shin Ramyun"
$ws.Range("O13").Value = "'0.484"
$ws.Range("P13").Value = "Sabor intenso, picante, umami, fideos gruesos, muy alto en sodio"
$ws.Range("Q13").Value = "Amy’s Macaroni & Cheese (frozen)"
$ws.Range("R13").Value = "'0.596"
$ws.Range("S13").Value = "Queso real, textura casera, sin conservadores, alto en grasa, algo caro"
$ws.Range("T13").Value = "Kraft Macaroni & Cheese Dinner"
$ws.Range("U13").Value = "'0.536"
$ws.Range("V13").Value = "Sabor nostálgico, clásico americano, fácil, no muy nutritivo, barato"
$ws.Range("W13").Value = "Velveeta Original Shells & Cheese (microwave cups)"
$ws.Range("X13").Value = "'0.520"
$ws.Range("Y13").Value = "Muy cremoso, porción individual, rápido, salado, ideal para niños"
$ws.Range("Z13").Value = "Wild Planet Wild Tuna Pasta Salad"
$ws.Range("AA13").Value = "'0.733"
$ws.Range("AB13").Value = "Sabor fresco, buena proteína, saludable, porción algo pequeña"
$ws.Range("AC13").Value = "StarKist Chicken Creations (Chicken Salad)"
$ws.Range("AD13").Value = "'0.594"
$ws.Range("AE13").Value = "Portátil, saludable, fácil, buena textura, sabor suave"
$ws.Range("AF13").Value = "Jack Link’s Beef Jerky Original"
$ws.Range("AG13").Value = "'0.570"
$ws.Range("AH13").Value = "Ahumado, sabroso, alto en proteína, snack ideal, porción pequeña"

# La celda G13 contiene JSON multilínea; Excel recalcula el alto de fila al
# escribirla. Lo fijamos de nuevo a la altura estándar de la hoja para que
# coincida con el resto de filas de datos.
$ws.Rows.Item(13).RowHeight = 15
